# error handling on scrapper_Service
#
# The scraper re-ran and picked up a corrected part number plus a new
# manufacturer reference row. Apply the same edits Excel would produce:
#   - A1: "3273112" -> "3273114" (kept as text, it is a part/catalog number)
#   - B1: the old leftover quantity (5) is cleared out
#   - A2 (new row): manufacturer reference "6SL32105BE211UV0", bold,
#     wrapped, with a taller row to fit the wrapped text
#   - selection left sitting on the newly entered cell A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: update the catalog number, keep it textual -------------------
# A plain numeric-looking string typed into a General cell is auto-typed
# as a number by Excel; prefixing with an apostrophe is the standard way
# to force it to stay text (quote-prefixed), matching the source data
# which stores this as a shared string.
$ws.Range("A1").Value = "'3273114"

# --- B1: clear the stale quantity value --------------------------------
$ws.Range("B1").ClearContents()

# --- A2: new manufacturer reference row --------------------------------
$ws.Range("A2").Value = "6SL32105BE211UV0"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").WrapText = $true

# Row heights: default row grew slightly, the new wrapped row is taller
$ws.Rows.Item(1).RowHeight = 15.65
$ws.Rows.Item(2).RowHeight = 44

# Leave the selection on the freshly-entered cell
[void]$ws.Range("A2").Select()
